$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the header row (row 1): copy it and insert the copy at row 2,
# pushing all existing data rows down by one.
$ws.Rows("1:1").Copy()
$ws.Rows("2:2").Insert()

# Update the view/selection to match the target state.
$ws.Range("P3").Select()
$excel.ActiveWindow.ScrollColumn = 4
